$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data set (header row 1 is unchanged).
# Columns: Posición, Nombre, Calificación General, Calificación Sommelier,
#          Precio Sugerido, Bodega, Región, País, Varietal
# NOTE: "Calificación Sommelier" and "Precio Sugerido" are stored as TEXT
# (e.g. "7.1", "5500.0") in the workbook, so they are entered with a
# leading apostrophe to force text instead of Excel's automatic number
# conversion.
$data = @(
    @(1,  "Vino El Cóndor",     7.0, "'7.1", "'5500.0", "Bodega Lopez",    "Zona Centro-Oeste", "Argentina", "Cabernet"),
    @(2,  "Gran Vino Mendoza",  7.0, "'4.7", "'3500.0", "Bodega Lopez",    "Zona Centro-Oeste", "Argentina", "Cabernet"),
    @(3,  "Vino Regional 2019", 7.0, "'4.3", "'2000.0", "Bodega Regional", "Zona Noreste",      "Argentina", "Malbec"),
    @(4,  "Gran Reserva 2015",  7.0, "'4.2", "'2500.0", "Bodega Lopez",    "Zona Centro-Oeste", "Argentina", "Malbec"),
    @(5,  "Vino Nacional",      7.0, "'4.1", "'1500.0", "Bodega Lopez",    "Zona Centro-Oeste", "Argentina", "Malbec"),
    @(6,  "San Juan Reserva",   7.0, "'4.1", "'3000.0", "Bodega Lopez",    "Zona Centro-Oeste", "Argentina", "Cabernet"),
    @(7,  "Tinto Gaucho",       7.0, "'3.5", "'1700.0", "Bodega Regional", "Zona Noreste",      "Argentina", "Cabernet"),
    @(8,  "Merlot Clásico",     7.0, "'3.1", "'1300.0", "Bodega Regional", "Zona Noreste",      "Argentina", "Cabernet"),
    @(9,  "Vino de La Rioja",   7.0, "'3.0", "'6000.0", "Bodega Regional", "Zona Noreste",      "Argentina", "Malbec"),
    @(10, "Vino Numero 8",      7.0, "'2.7", "'1500.0", "Bodega Lopez",    "Zona Centro-Oeste", "Argentina", "Malbec")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $ws.Cells.Item($rowIndex, 9).Value = $row[8]
    $rowIndex++
}
